$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "`'29.118.33"
$ws.Range("E2").Value = '  +0.42%  '

# Row 3
$ws.Range("D3").Value = "`'1.836.21"
$ws.Range("E3").Value = '  +0.36%  '

# Row 4
$ws.Range("D4").Value = "`'1.001"
$ws.Range("E4").Value = '  +0.37%  '

# Row 5
$ws.Range("D5").Value = "`'243.86"
$ws.Range("E5").Value = '  +0.02%  '

# Row 6
$ws.Range("E6").Value = '  -2.17%  '

# Row 7
$ws.Range("D7").Value = "`'1.003"
$ws.Range("E7").Value = '  +0.44%  '

# Row 8
$ws.Range("D8").Value = "`'0.07472"
$ws.Range("E8").Value = '  -0.56%  '

# Row 9
$ws.Range("D9").Value = "`'0.2920"
$ws.Range("E9").Value = '  -0.57%  '

# Row 10
$ws.Range("D10").Value = "`'23.13"
$ws.Range("E10").Value = '  +1.20%  '

# Row 11
$ws.Range("D11").Value = "`'0.07702"
$ws.Range("E11").Value = '  -0.28%  '

# Row 12
$ws.Range("D12").Value = "`'1.844.23"
$ws.Range("E12").Value = '  +0.05%  '

# Row 13
$ws.Range("D13").Value = "`'4.991"
$ws.Range("E13").Value = '  -0.03%  '

# Row 14
$ws.Range("D14").Value = "`'0.6711"
$ws.Range("E14").Value = '  +0.04%  '

# Row 15
$ws.Range("D15").Value = "`'82.59"
$ws.Range("E15").Value = '  -0.51%  '

# Row 16
$ws.Range("D16").Value = "`'0.000009319"
$ws.Range("E16").Value = '  -4.05%  '

# Row 17
$ws.Range("D17").Value = "`'5.933"
$ws.Range("E17").Value = '  -2.41%  '

# Row 18
$ws.Range("D18").Value = "`'29.124.04"
$ws.Range("E18").Value = '  +0.31%  '

# Row 19
$ws.Range("D19").Value = "`'2.105.22"
$ws.Range("E19").Value = '  +0.91%  '

# Row 20
$ws.Range("D20").Value = "`'231.22"
$ws.Range("E20").Value = '  +2.17%  '

# Row 21
$ws.Range("D21").Value = "`'12.63"
$ws.Range("E21").Value = '  +0.67%  '

# Row 22
$ws.Range("D22").Value = "`'1.004"
$ws.Range("E22").Value = '  +0.59%  '

# Row 23
$ws.Range("D23").Value = "`'7.165"
$ws.Range("E23").Value = '  -0.08%  '

# Row 24
$ws.Range("D24").Value = "`'1.003"
$ws.Range("E24").Value = '  +0.44%  '

# Row 25
$ws.Range("D25").Value = "`'160.16"
$ws.Range("E25").Value = '  +0.32%  '

# Row 26
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = "`'0.1387"
$ws.Range("E26").Value = '  -1.35%  '

# Row 27
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = "`'8.509"
$ws.Range("E27").Value = '  -0.25%  '

# Row 28
$ws.Range("D28").Value = "`'17.83"
$ws.Range("E28").Value = '  -0.38%  '

# Row 29
$ws.Range("D29").Value = "`'1.502"
$ws.Range("E29").Value = '  +0.24%  '

# Row 30
$ws.Range("D30").Value = "`'4.161"
$ws.Range("E30").Value = '  +1.04%  '

# Row 31
$ws.Range("D31").Value = "`'4.123"
$ws.Range("E31").Value = '  +1.33%  '

# Row 32
$ws.Range("D32").Value = "`'0.05526"
$ws.Range("E32").Value = '  +3.04%  '

# Row 33
$ws.Range("D33").Value = "`'1.201"
$ws.Range("E33").Value = '  +0.35%  '

# Row 34
$ws.Range("D34").Value = "`'0.7455"
$ws.Range("E34").Value = '  +0.30%  '

# Row 35
$ws.Range("D35").Value = "`'1.836"
$ws.Range("E35").Value = '  -1.17%  '

# Row 36
$ws.Range("D36").Value = "`'1.140"
$ws.Range("E36").Value = '  +0.17%  '

# Row 37
$ws.Range("D37").Value = "`'2.665"
$ws.Range("E37").Value = '  +0.46%  '

# Row 38
$ws.Range("D38").Value = "`'2.768"
$ws.Range("E38").Value = '  +0.55%  '

# Row 39
$ws.Range("D39").Value = "`'1.220.17"
$ws.Range("E39").Value = '  -1.95%  '

# Row 40
$ws.Range("D40").Value = "`'0.01779"
$ws.Range("E40").Value = '  -0.34%  '

# Row 41
$ws.Range("D41").Value = "`'6.462"
$ws.Range("E41").Value = '  -1.84%  '

# Row 42
$ws.Range("D42").Value = "`'0.8933"
$ws.Range("E42").Value = '  -0.85%  '

# Row 43
$ws.Range("D43").Value = "`'1.003"
$ws.Range("E43").Value = '  +0.40%  '

# Row 44
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = "`'101.84"
$ws.Range("E44").Value = '  +0.34%  '

# Row 45
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = "`'2.006.43"
$ws.Range("E45").Value = '  +1.03%  '

# Row 46
$ws.Range("D46").Value = "`'65.70"
$ws.Range("E46").Value = '  +1.48%  '

# Row 47
$ws.Range("D47").Value = "`'0.00000000121"
$ws.Range("E47").Value = '  -1.65%  '

# Row 48
$ws.Range("D48").Value = "`'0.5102"
$ws.Range("E48").Value = '  +0.02%  '

# Row 49
$ws.Range("D49").Value = "`'0.4071"
$ws.Range("E49").Value = '  +0.14%  '

# Row 50
$ws.Range("D50").Value = "`'9.094"
$ws.Range("E50").Value = '  +0.79%  '

# Row 51
$ws.Range("D51").Value = "`'0.05837"
$ws.Range("E51").Value = '  +1.26%  '
